# JD Unternehmen -> Veranstaltung & Kategorie geloescht
#
# The helper columns L:M (L1 header "Unternehmen", L2:L28 category ids,
# M2:M28 =COUNTIF(Wahlbereich,Lx) counts) are removed from the "Wahl"
# sheet. Removing them also drops the only reference to the shared
# string "Unternehmen", and shrinks the sheet's used range back down
# to A:I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wahl")

# Remove the whole helper block (header in row 1 + data rows 2-28)
$ws.Range("L1:M28").ClearContents()

# Match the selection Excel leaves behind after deleting the whole
# column(s) worth of content in L:M
$ws.Range("M1:M1048576").Select()
